$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price/Volume columns remain text so values like "1.00" or "0.140"
# keep their exact literal formatting instead of being parsed as numbers.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = '67.424.03'
$ws.Range("E2").Value = '  -0.14%  '
$ws.Range("D3").Value = '3.524.88'
$ws.Range("E3").Value = '  -0.16%  '
$ws.Range("E4").Value = '  +0.06%  '
$ws.Range("D5").Value = '613.64'
$ws.Range("E5").Value = '  +0.22%  '
$ws.Range("E6").Value = '  -1.74%  '
$ws.Range("D7").Value = '3.524.35'
$ws.Range("E7").Value = '  -0.10%  '
$ws.Range("D8").Value = '1.00'
$ws.Range("E8").Value = '  +0.03%  '
$ws.Range("D9").Value = '0.483'
$ws.Range("E9").Value = '  -0.61%  '
$ws.Range("E10").Value = '  -0.78%  '
$ws.Range("E11").Value = '  +3.47%  '
$ws.Range("D12").Value = '0.426'
$ws.Range("E12").Value = '  -1.38%  '
$ws.Range("E13").Value = '  -0.47%  '
$ws.Range("D14").Value = '32.11'
$ws.Range("E14").Value = '  -0.67%  '
$ws.Range("D15").Value = '4.120.87'
$ws.Range("E15").Value = '  -0.09%  '
$ws.Range("D16").Value = '3.537.03'
$ws.Range("E16").Value = '  +0.08%  '
$ws.Range("D17").Value = '67.418.77'
$ws.Range("E17").Value = '  -0.09%  '
$ws.Range("E18").Value = '  -0.11%  '
$ws.Range("D19").Value = '6.40'
$ws.Range("E19").Value = '  +0.41%  '
$ws.Range("D20").Value = '15.34'
$ws.Range("E20").Value = '  -1.56%  '
$ws.Range("D21").Value = '446.25'
$ws.Range("E21").Value = '  -1.93%  '
$ws.Range("D22").Value = '9.52'
$ws.Range("E22").Value = '  +1.19%  '
$ws.Range("D23").Value = '0.625'
$ws.Range("E23").Value = '  -2.87%  '
$ws.Range("D24").Value = '77.52'
$ws.Range("E24").Value = '  -1.52%  '
$ws.Range("E25").Value = '  +10.53%  '
$ws.Range("D26").Value = '3.665.39'
$ws.Range("E26").Value = '  -0.01%  '
$ws.Range("E27").Value = '  +0.02%  '
$ws.Range("D28").Value = '10.23'
$ws.Range("E28").Value = '  -2.86%  '
$ws.Range("E29").Value = '  +1.21%  '
$ws.Range("D30").Value = '2.51'
$ws.Range("E30").Value = '  -2.07%  '
$ws.Range("D31").Value = '1.57'
$ws.Range("E31").Value = '  -8.60%  '
$ws.Range("E32").Value = '  +0.07%  '
$ws.Range("E33").Value = '  +4.09%  '
$ws.Range("E34").Value = '  -0.70%  '
$ws.Range("E35").Value = '  -0.87%  '
$ws.Range("D36").Value = '3.517.04'
$ws.Range("E36").Value = '  -0.22%  '
$ws.Range("E37").Value = '  -3.62%  '
$ws.Range("E38").Value = '  -0.05%  '
$ws.Range("E39").Value = '  -0.03%  '
$ws.Range("E40").Value = '  +0.07%  '
$ws.Range("D41").Value = '176.86'
$ws.Range("E41").Value = '  +1.97%  '
$ws.Range("E42").Value = '  +1.75%  '
$ws.Range("D43").Value = '0.0884'
$ws.Range("E43").Value = '  +0.25%  '
$ws.Range("E44").Value = '  -3.52%  '
$ws.Range("E45").Value = '  -1.05%  '
$ws.Range("D46").Value = '28.19'
$ws.Range("E46").Value = '  -4.07%  '
$ws.Range("E47").Value = '  -1.37%  '
$ws.Range("D48").Value = '2.64'
$ws.Range("E48").Value = '  +1.09%  '
$ws.Range("E49").Value = '  +2.12%  '
$ws.Range("D50").Value = '7.60'
$ws.Range("D51").Value = '0.997'
$ws.Range("E51").Value = '  -2.41%  '
